# edit.ps1 - applies the diff to draft-ietf-pce-sr-bidir-path-02.pptx
#
# Changes:
#  1) Slide 4 ("PCE Initiated LSPs" diagram), shape "TextBox 26":
#       - reposition/resize the green-outlined callout box
#       - add a new bullet paragraph "B flag in RP object"
#  2) Slide 5 ("PCC Initiated LSPs" diagram), shape "TextBox 26":
#       - reposition/resize the green-outlined callout box (slightly
#         different target offset than slide 4)
#       - add the same new bullet paragraph "B flag in RP object"
#  3) Slide 8 ("Next Steps"), shape "Content Placeholder 2":
#       - remove the second paragraph "Add in Queue for WG LC?"
#
# NOTE: this COM-interop runtime's PowerShell engine does not bind named
# ("-Param value") arguments correctly for user-defined functions, so all
# helper calls below use positional arguments only.

# ---------------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height are exposed as single
# precision (float) values measured in points. Internally the host stores
# geometry in EMU (914400 EMU/inch, 12700 EMU/point) and truncates the
# point value (after it has been rounded to float32) when converting back
# to EMU. To land on an exact target EMU value we search nearby float32
# values for one that truncates to the desired EMU amount.
# ---------------------------------------------------------------------------
function Get-PointsForEmu($TargetEmu) {
    for ($delta = 0; $delta -le 400; $delta++) {
        foreach ($sign in 1, -1) {
            if ($delta -eq 0 -and $sign -eq -1) { continue }
            $candidatePt = ($TargetEmu / 12700.0) + ($sign * $delta * 0.0000005)
            $f32 = [float]$candidatePt
            $emu = [math]::Floor([double]$f32 * 12700.0)
            if ($emu -eq $TargetEmu) {
                return $candidatePt
            }
        }
    }
    return $TargetEmu / 12700.0
}

function Set-ShapeGeometryEmu($Shape, $OffX, $OffY, $ExtCx, $ExtCy) {
    $Shape.Left = Get-PointsForEmu $OffX
    $Shape.Top = Get-PointsForEmu $OffY
    $Shape.Width = Get-PointsForEmu $ExtCx
    $Shape.Height = Get-PointsForEmu $ExtCy
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 4 - TextBox 26
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$box4 = $slide4.Shapes.Item("TextBox 26")

Set-ShapeGeometryEmu $box4 2795328 2419350 3400764 954107

$tr4 = $box4.TextFrame.TextRange
$null = $tr4.InsertAfter([char]13 + "B flag in RP object")

# ---------------------------------------------------------------------------
# 2) Slide 5 - TextBox 26
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$box5 = $slide5.Shapes.Item("TextBox 26")

Set-ShapeGeometryEmu $box5 2798703 2434805 3400764 954107

$tr5 = $box5.TextFrame.TextRange
$null = $tr5.InsertAfter([char]13 + "B flag in RP object")

# ---------------------------------------------------------------------------
# 3) Slide 8 - Content Placeholder 2: drop "Add in Queue for WG LC?" para
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$content8 = $slide8.Shapes.Item("Content Placeholder 2")
$content8.TextFrame.TextRange.Text = "Welcome your review comments and suggestions"

Write-Host "edit.ps1 applied successfully"
